# "15 - Strings.pptx", slide 8 ("Grammar Rules Relevant to Strings").
#
# In the "Content Placeholder 2" body text, the grammar rule for
# `variable` reads:
#     variable  = ( varId | paramId) { indexExpr | fieldExpr } .
# A space is inserted right before the closing ")" so it reads:
#     variable  = ( varId | paramId ) { indexExpr | fieldExpr } .

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(8)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# Locate "paramId) " within the shape's text so we don't depend on a
# hard-coded character offset.
$fullText = $tr.Text
$needle   = "paramId) "
$pos0     = $fullText.IndexOf($needle)

if ($pos0 -ge 0) {
    # 1-based index (PowerPoint TextRange.Characters is 1-based) of the
    # ")" immediately after "paramId", plus the space right after it.
    $start = $pos0 + "paramId".Length + 1
    $closeParenAndSpace = $tr.Characters($start, 2)
    $closeParenAndSpace.Text = " ) "
}
